$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the shared string for the "ඕ" character: it previously had a stray
# leading tab character (shared with the adjacent "ඔ" entry's formatting).
# Setting the clean value causes the workbook to drop the old shared-string
# entry and append the corrected one at the end of the table, which is what
# the target diff shows (old index 10 removed, new entry appended as index 30,
# and every row referencing a shifted shared string is renumbered accordingly).
$ws.Range("B11").Value = "ඕ"

# Match the final selection left behind in the saved workbook.
$ws.Range("G8").Select()
